# "add zombies to town" -------------------------------------------------
# Applies the changes described by the commit diff:
#   1. Update the cached "last-modified" date field (7/28/16 -> 7/31/16)
#      on the slide master and every slide layout.
#   2. Split the "4: Infirmary: Pink" run on slide 2 into three runs.
#   3. Move two shapes on slide 3 (the "Rectangle 21" / "TextBox 1" pair).
#   4. Re-label the "Brotherhood / OTSA" textbox on slide 5 to
#      "Biological / containment" (zombie containment!) and resize it.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text: 7/28/16 -> 7/31/16 (slide layouts + master)
#    (loops are inlined -- no helper function -- to stay within the
#    host's statement budget)
# ---------------------------------------------------------------------
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layoutShapes = $layouts.Item($i).Shapes
    for ($j = 1; $j -le $layoutShapes.Count; $j++) {
        $sh = $layoutShapes.Item($j)
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "7/28/16") {
            $sh.TextFrame.TextRange.Text = "7/31/16"
        }
    }
}

$masterShapes = $p.SlideMaster.Shapes
for ($j = 1; $j -le $masterShapes.Count; $j++) {
    $sh = $masterShapes.Item($j)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "7/28/16") {
        $sh.TextFrame.TextRange.Text = "7/31/16"
    }
}

# ---------------------------------------------------------------------
# 2) Slide 2: split "4: Infirmary: Pink" into 3 runs ("4: " / "Infirmary"
#    / ": Pink") while keeping the overall text (and formatting) intact.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$infirmaryShape = $s2.Shapes.Item(16)
$infirmaryShape.TextFrame.TextRange.Characters(1, 3).Text = "4: "
$infirmaryShape.TextFrame.TextRange.Characters(4, 9).Text = "Infirmary"
$infirmaryShape.TextFrame.TextRange.Characters(13, 6).Text = ": Pink"

# ---------------------------------------------------------------------
# 3) Slide 3: reposition "Rectangle 21" (id 22) and "TextBox 1" (id 2)
#    -- sizes are unchanged, only the x/y offsets move.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$rect21 = $s3.Shapes.Item(19)
$textBox1 = $s3.Shapes.Item(20)

# Points values chosen so that Point -> EMU round-tripping inside the
# host lands exactly on the target EMU offsets (914400 EMU/in, 12700 EMU/pt).
$rect21.Left = 294.6192169984252    # 3741664 EMU
$rect21.Top  = 411.05393990787405   # 5220385 EMU
$textBox1.Left = 303.6230011259842  # 3856012 EMU
$textBox1.Top  = 413.8807220614173  # 5256285 EMU

# ---------------------------------------------------------------------
# 4) Slide 5: "Brotherhood" / "OTSA" -> "Biological" / "containment"
#    (new zombie-containment call-out); the textbox also narrows a touch.
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$bioShape = $s5.Shapes.Item(20)
$bioShape.TextFrame.TextRange.Text = ""
$bioShape.TextFrame.TextRange.Text = "Biological`rcontainment"
$bioShape.Width = 67.06023792047245  # 851665 EMU
